$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of the existing header cell H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the data values for the new I and J columns (rows 2-12)
$iValues = @(7, 1, 1, 1, 1, 1, 1, 1, 1, 5, 5)
$jValues = @(8, 2, 6, 5, 6, 6, 4, 5, 4, 6, 6)

for ($r = 2; $r -le 12; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
